$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.899.09'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '1.894.95'
$ws.Range("E3").Value = '  +1.47%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4790'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2915'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.46%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06587'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.48'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07796'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '97.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.24%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.897.31'
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7413'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +7.48%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.203'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '282.75'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.31%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '30.881.11'
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007646'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.154.18'
$ws.Range("E22").Value = '  +1.45%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.325'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9995'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.265'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.89%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.393'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.80%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.70%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.997'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.50%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.379'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1005'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.522'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.63%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.405'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.52%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.151'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.95%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.04791'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.136'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7092'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.718'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01882'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.776'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.491'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.68%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.945'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4237'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.13%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8498'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.17%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9996'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.518'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.12%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.186'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.68%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '953.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.96%  '
